$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 37 and 38 swap content (Monero <-> ImmutableX) with updated price/volume
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("E37").Value = "  +6.08%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "167.77"
$ws.Range("E38").Value = "  +1.32%  "

# Remaining price/volume updates
$ws.Range("D2").Value = "62.143.36"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "3.427.84"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.10"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.12"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.62"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "4.015.19"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.81"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "3.430.41"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "62.224.23"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.06"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.22"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.05"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.90"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "3.568.12"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.68"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.02"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D39").Value = "3.461.17"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "28.91"
$ws.Range("E40").Value = "  +9.42%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").Value = "2.516.53"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.12"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("E51").Value = "  +0.89%  "
